$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Summary block: total "Valor Mora" and "Cant. Trabajadores" counters
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 568229
$ws.Range("C13").Value = 8

# ---------------------------------------------------------------------------
# 2) Expand the worker table from 6 rows (16-21) to 9 rows (16-24).
#    Insert 3 new rows right after the existing last data row (21) so the
#    trailing signature block (previously rows 26-27) is pushed down to
#    rows 29-30, matching the target layout.
# ---------------------------------------------------------------------------
$ws.Rows("22:24").Insert()

# Row 21 still carries the "last row" (thicker bottom border) formatting.
# Re-home that formatting on the new last row (24).
$ws.Range("B21:J21").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)

# Row 20 carries the regular interior-row formatting; stamp it onto rows
# 21-23 (row 21 needs to flip from "last row" style back to a normal one
# now that row 24 is the new last row).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B22:J22").PasteSpecial(-4122)
$ws.Range("B23:J23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Replace the table contents (rows 16-24) with the refreshed data set.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora,
#    G=Salario Basico
# ---------------------------------------------------------------------------
$data = @(
  @("CC", "3830222",    "ELKIN ENRIQUE ARRIETA CAICEDO",    "2507",  86000, 2150000),
  @("CC", "1143374161", "EDWIN ALVAREZ VASQUEZ",            "2507",  62910, 1572740),
  @("CC", "1043636663", "NEYDER JOSE CAMERA ULLOA",         "2507",  43853, 1495000),
  @("CC", "1047415641", "HENRY AMALIO CASTILLO TORRES",     "2312",  11200, 1400000),
  @("CC", "1047415641", "HENRY AMALIO CASTILLO TORRES",     "2311",  20533, 1400000),
  @("CC", "73149881",   "HECTOR WILLIAM PORRAS BARBOZA",    "2507", 148000, 3700000),
  @("CC", "1143348470", "GABRIEL ANGEL PEREZ RONDANO",      "2507", 116000, 2900000),
  @("CC", "1001971665", "DEYMER JOSE SANMARTIN GUTIERREZ",  "2507",  59800, 1495000),
  @("CC", "1002059655", "DEIMER DE JESUS GARIZAO GARIZAO",  "2507",  19933, 1495000)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
  $ws.Cells.Item($r, 7).Value = $row[5]
  $r = $r + 1
}
